$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 183.4
$ws.Range("I9").Value = 100
$ws.Range("J9").Value = 266.8
$ws.Range("K9").Value = 100
$ws.Range("L9").Value = 266.8
$ws.Range("M9").Value = 69
$ws.Range("N9").Value = -604.8
$ws.Range("H80").Value = 2217.5527
$ws.Range("I80").Value = 880.9048
$ws.Range("K80").Value = 2642.7144
$ws.Range("M80").Value = -1644.7144
$ws.Range("H83").Value = 2217.5527
$ws.Range("I83").Value = 880.9048
$ws.Range("K83").Value = 7928.1432
$ws.Range("M83").Value = -2936.1432
$ws.Range("H132").Value = 3525.795
$ws.Range("I132").Value = 3578.2354
$ws.Range("K132").Value = 10734.7062
$ws.Range("M132").Value = -8204.706200000001
$ws.Range("H137").Value = 1929.091
$ws.Range("I137").Value = 2047.5883
$ws.Range("K137").Value = 6142.7649
$ws.Range("M137").Value = -3592.7649

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H75").Value = 27000
$ws.Range("J75").Value = 27000
$ws.Range("L75").Value = 27000
$ws.Range("N75").Value = -28748
$ws.Range("H78").Value = 27000
$ws.Range("J78").Value = 27000
$ws.Range("L78").Value = 81000
$ws.Range("N78").Value = -89736
$ws.Range("H102").Value = 4072.9092
$ws.Range("I102").Value = 2002.5
$ws.Range("J102").Value = 5256
$ws.Range("K102").Value = 2002.5
$ws.Range("L102").Value = 5256
$ws.Range("M102").Value = -380.5
$ws.Range("N102").Value = -8500
$ws.Range("H110").Value = 915.5333000000001
$ws.Range("I110").Value = 915.5333000000001
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 915.5333000000001
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1129.4667
$ws.Range("N110").ClearContents()
$ws.Range("H132").Value = 18538122
$ws.Range("I132").Value = 25641634
$ws.Range("K132").Value = 76924902
$ws.Range("M132").Value = -76922372

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2780.0454
$ws.Range("I20").Value = 2417.4285
$ws.Range("J20").Value = 3414.625
$ws.Range("K20").Value = 2417.4285
$ws.Range("L20").Value = 3414.625
$ws.Range("M20").Value = -2170.4285
$ws.Range("N20").Value = -3908.625
$ws.Range("H82").Value = 15452
$ws.Range("I82").Value = 7459.2
$ws.Range("J82").Value = 21161.143
$ws.Range("K82").Value = 7459.2
$ws.Range("L82").Value = 21161.143
$ws.Range("M82").Value = -7076.2
$ws.Range("N82").Value = -21927.143
$ws.Range("H85").Value = 15452
$ws.Range("I85").Value = 7459.2
$ws.Range("J85").Value = 21161.143
$ws.Range("K85").Value = 7459.2
$ws.Range("L85").Value = 21161.143
$ws.Range("M85").Value = -6133.2
$ws.Range("N85").Value = -23813.143
$ws.Range("H86").Value = 45459964
$ws.Range("I86").Value = 71431150
$ws.Range("J86").Value = 10376.75
$ws.Range("K86").Value = 71431150
$ws.Range("L86").Value = 10376.75
$ws.Range("M86").Value = -71430027
$ws.Range("N86").Value = -12622.75
$ws.Range("H89").Value = 45459964
$ws.Range("I89").Value = 71431150
$ws.Range("J89").Value = 10376.75
$ws.Range("K89").Value = 357155750
$ws.Range("L89").Value = 51883.75
$ws.Range("M89").Value = -357150134
$ws.Range("N89").Value = -63115.75
$ws.Range("H134").Value = 7670.8276
$ws.Range("I134").Value = 2746.9583
$ws.Range("K134").Value = 8240.874899999999
$ws.Range("M134").Value = -5705.874899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9972.727999999999
$ws.Range("J4").Value = 9972.727999999999
$ws.Range("L4").Value = 9972.727999999999
$ws.Range("N4").Value = -10196.728
$ws.Range("H16").Value = 1389.9546
$ws.Range("I16").Value = 1331.5834
$ws.Range("J16").Value = 1460
$ws.Range("K16").Value = 1331.5834
$ws.Range("L16").Value = 1460
$ws.Range("M16").Value = -1044.5834
$ws.Range("N16").Value = -2034
$ws.Range("H58").Value = 859.8148
$ws.Range("I58").Value = 801
$ws.Range("J58").Value = 1065.6666
$ws.Range("K58").Value = 801
$ws.Range("L58").Value = 1065.6666
$ws.Range("M58").Value = -598
$ws.Range("N58").Value = -1471.6666
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H105").Value = 852.8570999999999
$ws.Range("I105").Value = 867.2727
$ws.Range("J105").Value = 800
$ws.Range("K105").Value = 867.2727
$ws.Range("L105").Value = 800
$ws.Range("M105").Value = 879.7273
$ws.Range("N105").Value = -4294
$ws.Range("H113").Value = 1389.9546
$ws.Range("I113").Value = 1331.5834
$ws.Range("J113").Value = 1460
$ws.Range("K113").Value = 1331.5834
$ws.Range("L113").Value = 1460
$ws.Range("M113").Value = 838.4166
$ws.Range("N113").Value = -5800
$ws.Range("H136").Value = 859.8148
$ws.Range("I136").Value = 801
$ws.Range("J136").Value = 1065.6666
$ws.Range("K136").Value = 2403
$ws.Range("L136").Value = 3196.9998
$ws.Range("M136").Value = 147
$ws.Range("N136").Value = -8296.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 449.57144
$ws.Range("I107").Value = 366
$ws.Range("J107").Value = 496
$ws.Range("K107").Value = 1098
$ws.Range("L107").Value = 1488
$ws.Range("M107").Value = 822
$ws.Range("N107").Value = -5328

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4624.82
$ws.Range("I132").Value = 3298.4666
$ws.Range("K132").Value = 9895.399800000001
$ws.Range("M132").Value = -7365.399800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H16").Value = 406.92856
$ws.Range("I16").Value = 377
$ws.Range("K16").Value = 377
$ws.Range("M16").Value = -207
$ws.Range("H60").Value = 6000
$ws.Range("J60").Value = 6000
$ws.Range("L60").Value = 6000
$ws.Range("N60").Value = -7018
$ws.Range("H132").Value = 19571.725
$ws.Range("I132").Value = 29775.432
$ws.Range("J132").Value = 1593.762
$ws.Range("K132").Value = 89326.296
$ws.Range("L132").Value = 4781.286
$ws.Range("M132").Value = -86796.296
$ws.Range("N132").Value = -9841.286
$ws.Range("H136").Value = 7258.1665
$ws.Range("J136").Value = 6724
$ws.Range("L136").Value = 20172
$ws.Range("N136").Value = -25272

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 28720500
$ws.Range("I132").Value = 38182920
$ws.Range("J132").Value = 2698846.8
$ws.Range("K132").Value = 114548760
$ws.Range("L132").Value = 8096540.399999999
$ws.Range("M132").Value = -114546230
$ws.Range("N132").Value = -8101600.399999999
